# Bitacora_tareas.xlsx - "Errores de clasificación: solapa, edición y flags"
# Adds three new log rows to the "Log" sheet and two new summary rows to the
# "Resumen" sheet, documenting the new error-classification tab, inline
# record editing, and the editado/editado_detalle fields.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Log": append rows 22-24
# ---------------------------------------------------------------------------
$log = $wb.Worksheets.Item("Log")

$log.Range("A22").Value = "27/02/2025"
$log.Range("B22").Value = "16:50"
$log.Range("C22").Value = "Recategorización Alquiler → Alquileres y Servicios"
$log.Range("D22").Value = "Si la categoría original es exactamente Alquiler, el dashboard la muestra como Alquileres y Servicios (solo cambio de etiqueta visual, los números y agrupaciones siguen conciliando)."
$log.Range("E22").Value = "Diagnostico"

$log.Range("A23").Value = "27/02/2025"
$log.Range("B23").Value = "17:00"
$log.Range("C23").Value = "Detección de errores de clasificación (Egresos)"
$log.Range("D23").Value = 'Para egresos: si la descripción (más cat_desc/cliente) no contiene palabras relevantes de la categoría mostrada o de la cuenta contable, se recategoriza visualmente como Sin categoría y se registra como error de tipo "Inconsistencia entre Categoria , Cuenta Contable y Descripcion". En el modal mensual se agrega solapa Errores con el conteo y un acceso a un modal de detalle con todos los registros en error.'
$log.Range("E23").Value = "Diagnostico"

$log.Range("A24").Value = "27/02/2025"
$log.Range("B24").Value = "17:30"
$log.Range("C24").Value = "Modal errores: ampliar, editar registro y campos editado/editado_detalle"
$log.Range("D24").Value = "Ampliar modal de detalle de errores. Agregar icono de edición por registro que abre modal para actualizar en BD: Categoría y Cuenta contable solo desde valores existentes (dropdown), Descripción libre. Tabla transacciones: nuevos campos editado (flag) y editado_detalle (ej. Categoria, Descripcion, Cuenta Contable). Migración supabase_transacciones_editado.sql. Export Excel incluye editado y editado_detalle."
$log.Range("E24").Value = "Diagnostico"

# ---------------------------------------------------------------------------
# Sheet "Resumen": append rows 29-30
# ---------------------------------------------------------------------------
$resumen = $wb.Worksheets.Item("Resumen")

$resumen.Range("A29").Value = "Edición desde modal Errores"
$resumen.Range("B29").Value = "En el detalle de errores, icono de edición por registro. Abre modal para corregir: Categoría y Cuenta contable solo desde valores existentes en BD; Descripción libre. Al guardar se actualiza la fila y se marcan editado y editado_detalle (qué campos se editaron)."

$resumen.Range("A30").Value = "Campos editado y editado_detalle"
$resumen.Range("B30").Value = 'En transacciones: editado (boolean) y editado_detalle (texto, ej. "Categoria, Descripcion, Cuenta Contable"). Migración supabase_transacciones_editado.sql. Export Excel los incluye.'
